$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 7153.6665
$ws.Range("I41").Value = 183.33333
$ws.Range("J41").Value = 10638.833
$ws.Range("K41").Value = 183.33333
$ws.Range("L41").Value = 10638.833
$ws.Range("M41").Value = 256.66667
$ws.Range("N41").Value = -11518.833
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H106").Value = 38539580
$ws.Range("I106").Value = 92138.17999999999
$ws.Range("J106").Value = 250000500
$ws.Range("K106").Value = 92138.17999999999
$ws.Range("L106").Value = 250000500
$ws.Range("M106").Value = -91507.17999999999
$ws.Range("N106").Value = -250001762
$ws.Range("H113").Value = 2651.25
$ws.Range("I113").Value = 2363.077
$ws.Range("J113").Value = 3186.4285
$ws.Range("K113").Value = 2363.077
$ws.Range("L113").Value = 3186.4285
$ws.Range("M113").Value = 890.9229999999998
$ws.Range("N113").Value = -9694.4285
$ws.Range("H116").Value = 1429.1428
$ws.Range("I116").Value = 1429.1428
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1429.1428
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 2012.8572
$ws.Range("N116").ClearContents()
$ws.Range("H131").Value = 1552.1428
$ws.Range("I131").Value = 916.25
$ws.Range("J131").Value = 2400
$ws.Range("K131").Value = 2748.75
$ws.Range("L131").Value = 7200
$ws.Range("M131").Value = 2291.25
$ws.Range("N131").Value = -17280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 481.0909
$ws.Range("I97").Value = 346.66666
$ws.Range("J97").Value = 642.4
$ws.Range("K97").Value = 346.66666
$ws.Range("L97").Value = 642.4
$ws.Range("M97").Value = 149.33334
$ws.Range("N97").Value = -1634.4
$ws.Range("H102").Value = 1342.9259
$ws.Range("I102").Value = 1466.3182
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 1466.3182
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = 155.6818000000001
$ws.Range("N102").Value = -4044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 705.6
$ws.Range("I94").Value = 651.2727
$ws.Range("K94").Value = 651.2727
$ws.Range("M94").Value = -200.2727
$ws.Range("H99").Value = 635.26086
$ws.Range("I99").Value = 450
$ws.Range("K99").Value = 450
$ws.Range("M99").Value = 1048

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7694338.5
$ws.Range("I31").Value = 1889.4
$ws.Range("J31").Value = 33335836
$ws.Range("K31").Value = 1889.4
$ws.Range("L31").Value = 33335836
$ws.Range("M31").Value = -1594.4
$ws.Range("N31").Value = -33336426
$ws.Range("H34").Value = 7694338.5
$ws.Range("I34").Value = 1889.4
$ws.Range("J34").Value = 33335836
$ws.Range("K34").Value = 1889.4
$ws.Range("L34").Value = 33335836
$ws.Range("M34").Value = -1687.4
$ws.Range("N34").Value = -33336240
$ws.Range("H58").Value = 1353.4667
$ws.Range("I58").Value = 1235.8572
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 1235.8572
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -1032.8572
$ws.Range("N58").Value = -3406
$ws.Range("H111").Value = 29909.092
$ws.Range("J111").Value = 29909.092
$ws.Range("L111").Value = 29909.092
$ws.Range("N111").Value = -38089.092
$ws.Range("H132").Value = 2693.85
$ws.Range("I132").Value = 1160
$ws.Range("J132").Value = 5542.4287
$ws.Range("K132").Value = 3480
$ws.Range("L132").Value = 16627.2861
$ws.Range("M132").Value = -950
$ws.Range("N132").Value = -21687.2861
$ws.Range("H136").Value = 1353.4667
$ws.Range("I136").Value = 1235.8572
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3707.5716
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1157.5716
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 595.7308
$ws.Range("I5").Value = 321
$ws.Range("J5").Value = 831.2143
$ws.Range("K5").Value = 963
$ws.Range("L5").Value = 2493.6429
$ws.Range("M5").Value = -851
$ws.Range("N5").Value = -2717.6429
$ws.Range("H131").Value = 7027.0586
$ws.Range("J131").Value = 850
$ws.Range("L131").Value = 2550
$ws.Range("N131").Value = -12630
$ws.Range("H135").Value = 595.7308
$ws.Range("I135").Value = 321
$ws.Range("J135").Value = 831.2143
$ws.Range("K135").Value = 2889
$ws.Range("L135").Value = 7480.928699999999
$ws.Range("M135").Value = -354
$ws.Range("N135").Value = -12550.9287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 37171.668
$ws.Range("I34").Value = 10000
$ws.Range("K34").Value = 10000
$ws.Range("M34").Value = -9732
$ws.Range("H70").Value = 18218230
$ws.Range("I70").Value = 19619266
$ws.Range("J70").Value = 4750
$ws.Range("K70").Value = 19619266
$ws.Range("L70").Value = 4750
$ws.Range("M70").Value = -19618996
$ws.Range("N70").Value = -5290
$ws.Range("H73").Value = 18218230
$ws.Range("I73").Value = 19619266
$ws.Range("J73").Value = 4750
$ws.Range("K73").Value = 19619266
$ws.Range("L73").Value = 4750
$ws.Range("M73").Value = -19618330
$ws.Range("N73").Value = -6622
$ws.Range("H76").Value = 37171.668
$ws.Range("I76").Value = 10000
$ws.Range("K76").Value = 10000
$ws.Range("M76").Value = -9685
$ws.Range("H79").Value = 37171.668
$ws.Range("I79").Value = 10000
$ws.Range("K79").Value = 10000
$ws.Range("M79").Value = -8908
$ws.Range("H97").Value = 769.9545000000001
$ws.Range("J97").Value = 802.375
$ws.Range("L97").Value = 802.375
$ws.Range("N97").Value = -1794.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H132").Value = 1917.4546
$ws.Range("I132").Value = 1765.3334
$ws.Range("J132").Value = 2323.111
$ws.Range("K132").Value = 5296.0002
$ws.Range("L132").Value = 6969.333
$ws.Range("M132").Value = -2766.0002
$ws.Range("N132").Value = -12029.333
$ws.Range("H136").Value = 4571.724
$ws.Range("I136").Value = 4994.7827
$ws.Range("J136").Value = 2950
$ws.Range("K136").Value = 14984.3481
$ws.Range("L136").Value = 8850
$ws.Range("M136").Value = -12434.3481
$ws.Range("N136").Value = -13950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -31240
$ws.Range("H132").Value = 1714.1786
$ws.Range("I132").Value = 1345.7949
$ws.Range("J132").Value = 2559.2942
$ws.Range("K132").Value = 4037.384700000001
$ws.Range("L132").Value = 7677.882599999999
$ws.Range("M132").Value = -1507.384700000001
$ws.Range("N132").Value = -12737.8826
$ws.Range("H136").Value = 8439.177
$ws.Range("I136").Value = 8439.177
$ws.Range("K136").Value = 25317.531
$ws.Range("M136").Value = -22767.531

Write-Output "Applied Garuda_Profits price refresh edits"